{"js": "// Replace the date line and every \"AxB=\" multiplication problem text\n// with the updated values from the commit, via search-and-replace so\n// each run's existing formatting (font, size, etc.) is preserved.\nconst replacements = [\n  [\"2024-03-16 Saturday\", \"2024-03-17 Sunday\"],\n  [\"558\u00d79=\", \"171\u00d79=\"],\n  [\"631\u00d76=\", \"982\u00d79=\"],\n  [\"216\u00d74=\", \"754\u00d72=\"],\n  [\"160\u00d74=\", \"949\u00d75=\"],\n  [\"538\u00d74=\", \"821\u00d72=\"],\n  [\"388\u00d77=\", \"375\u00d75=\"],\n  [\"274\u00d72=\", \"799\u00d77=\"],\n  [\"298\u00d77=\", \"628\u00d74=\"],\n  [\"417\u00d72=\", \"980\u00d78=\"],\n  [\"334\u00d76=\", \"323\u00d78=\"],\n  [\"394\u00d76=\", \"543\u00d72=\"],\n  [\"820\u00d72=\", \"949\u00d73=\"],\n  [\"217\u00d77=\", \"422\u00d73=\"],\n  [\"337\u00d79=\", \"851\u00d77=\"],\n  [\"959\u00d75=\", \"406\u00d73=\"],\n  [\"559\u00d76=\", \"894\u00d78=\"],\n  [\"728\u00d72=\", \"405\u00d79=\"],\n  [\"649\u00d73=\", \"736\u00d79=\"],\n  [\"289\u00d75=\", \"191\u00d75=\"],\n  [\"920\u00d79=\", \"319\u00d78=\"],\n  [\"650\u00d74=\", \"778\u00d76=\"],\n  [\"111\u00d78=\", \"315\u00d76=\"],\n  [\"686\u00d74=\", \"898\u00d73=\"],\n  [\"461\u00d78=\", \"711\u00d76=\"],\n  [\"995\u00d76=\", \"768\u00d76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and every \"AxB=\" multiplication problem to the\n# new values from the commit, using Find/Replace so the existing run\n# formatting (font, size, etc.) on each piece of text is preserved.\n$d = $word.ActiveDocument\n\nfunction Replace-Text($oldText, $newText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($null, $false, $false, $false, $null, $null, $true, $null, $null, $newText, 2)\n}\n\nReplace-Text \"2024-03-16 Saturday\" \"2024-03-17 Sunday\"\n\nReplace-Text \"558\u00d79=\" \"171\u00d79=\"\nReplace-Text \"631\u00d76=\" \"982\u00d79=\"\nReplace-Text \"216\u00d74=\" \"754\u00d72=\"\nReplace-Text \"160\u00d74=\" \"949\u00d75=\"\nReplace-Text \"538\u00d74=\" \"821\u00d72=\"\n\nReplace-Text \"388\u00d77=\" \"375\u00d75=\"\nReplace-Text \"274\u00d72=\" \"799\u00d77=\"\nReplace-Text \"298\u00d77=\" \"628\u00d74=\"\nReplace-Text \"417\u00d72=\" \"980\u00d78=\"\nReplace-Text \"334\u00d76=\" \"323\u00d78=\"\n\nReplace-Text \"394\u00d76=\" \"543\u00d72=\"\nReplace-Text \"820\u00d72=\" \"949\u00d73=\"\nReplace-Text \"217\u00d77=\" \"422\u00d73=\"\nReplace-Text \"337\u00d79=\" \"851\u00d77=\"\nReplace-Text \"959\u00d75=\" \"406\u00d73=\"\n\nReplace-Text \"559\u00d76=\" \"894\u00d78=\"\nReplace-Text \"728\u00d72=\" \"405\u00d79=\"\nReplace-Text \"649\u00d73=\" \"736\u00d79=\"\nReplace-Text \"289\u00d75=\" \"191\u00d75=\"\nReplace-Text \"920\u00d79=\" \"319\u00d78=\"\n\nReplace-Text \"650\u00d74=\" \"778\u00d76=\"\nReplace-Text \"111\u00d78=\" \"315\u00d76=\"\nReplace-Text \"686\u00d74=\" \"898\u00d73=\"\nReplace-Text \"461\u00d78=\" \"711\u00d76=\"\nReplace-Text \"995\u00d76=\" \"768\u00d76=\"\n"}
